# Insert 3 new weekly price rows for "Feria Lagunitas de Puerto Montt - Palta"
# above the existing row 439 (shifting the old rows 439-523 down to 442-526),
# then populate the 3 new rows with the new survey data (fecha 44782, origin Peru).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 439; everything below shifts down.
$ws.Rows("439:441").Insert()

# --- New row 439: Especial / Peru / $/bandeja 10 kilos ---
$ws.Range("A439").Value = 4
$ws.Range("B439").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C439").Value = "Los Lagos"
$ws.Range("D439").Value = 44782
$ws.Range("E439").Value = 10
$ws.Range("F439").Value = "Fruta"
$ws.Range("G439").Value = 100106
$ws.Range("H439").Value = "Oleaginosos"
$ws.Range("I439").Value = 100106002
$ws.Range("J439").Value = "Palta"
$ws.Range("K439").Value = "Hass"
$ws.Range("L439").Value = "Especial"
$ws.Range("M439").Value = 150
$ws.Range("N439").Value = 35000
$ws.Range("O439").Value = 35000
$ws.Range("P439").Value = 35000
$ws.Range("Q439").Value = "`$/bandeja 10 kilos"
$ws.Range("R439").Value = "Perú"
$ws.Range("S439").Value = 3500
$ws.Range("T439").Value = 10

# --- New row 440: Primera / Peru / $/bandeja 10 kilos ---
$ws.Range("A440").Value = 4
$ws.Range("B440").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C440").Value = "Los Lagos"
$ws.Range("D440").Value = 44782
$ws.Range("E440").Value = 10
$ws.Range("F440").Value = "Fruta"
$ws.Range("G440").Value = 100106
$ws.Range("H440").Value = "Oleaginosos"
$ws.Range("I440").Value = 100106002
$ws.Range("J440").Value = "Palta"
$ws.Range("K440").Value = "Hass"
$ws.Range("L440").Value = "Primera"
$ws.Range("M440").Value = 150
$ws.Range("N440").Value = 33000
$ws.Range("O440").Value = 33000
$ws.Range("P440").Value = 33000
$ws.Range("Q440").Value = "`$/bandeja 10 kilos"
$ws.Range("R440").Value = "Perú"
$ws.Range("S440").Value = 3300
$ws.Range("T440").Value = 10

# --- New row 441: Segunda / Peru / $/bandeja 10 kilos ---
$ws.Range("A441").Value = 4
$ws.Range("B441").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C441").Value = "Los Lagos"
$ws.Range("D441").Value = 44782
$ws.Range("E441").Value = 10
$ws.Range("F441").Value = "Fruta"
$ws.Range("G441").Value = 100106
$ws.Range("H441").Value = "Oleaginosos"
$ws.Range("I441").Value = 100106002
$ws.Range("J441").Value = "Palta"
$ws.Range("K441").Value = "Hass"
$ws.Range("L441").Value = "Segunda"
$ws.Range("M441").Value = 150
$ws.Range("N441").Value = 30000
$ws.Range("O441").Value = 30000
$ws.Range("P441").Value = 30000
$ws.Range("Q441").Value = "`$/bandeja 10 kilos"
$ws.Range("R441").Value = "Perú"
$ws.Range("S441").Value = 3000
$ws.Range("T441").Value = 10
